$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching style of existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I and J column values for rows 2-48
$iValues = @(11,10,9,11,8,7,8,9,8,7,7,8,9,7,5,7,7,7,6,7,7,8,7,8,7,6,8,6,8,7,8,7,6,6,9,9,8,8,5,7,8,7,8,1,7,4,3)
$jValues = @(11,10,9,11,8,7,8,9,8,8,8,8,9,7,6,7,7,7,6,7,7,8,7,8,7,7,8,7,8,7,8,7,7,6,9,10,8,8,7,8,8,7,8,4,9,6,4)

for ($r = 2; $r -le 48; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
